$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6 for the new "Tree Canopy Cover" (USGS/GEE) entry.
# Everything currently on row 6+ (Wildfire Risk, Flood Risk, Ozone, ...) shifts down by one.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6.
$ws.Range("A6").Value = "Tree Canopy Cover"
$ws.Range("C6").Value = "Y"
$ws.Range("D6").Value = "Earth Engine"
$ws.Range("E6").Value = "30m"
$ws.Range("G6").Value = "USGS/GEE"
$ws.Range("H6").Value = "https://developers.google.com/earth-engine/datasets/catalog/USGS_NLCD_RELEASES_2016_REL"
$ws.Range("I6").Value = "AK and HI are separate images, band for percent tree canopy cover. Code example of importing and mapping: https://code.earthengine.google.com/0530abe36a81b15ce425b3a11c81cefe "

# Update the note on the Flood Risk row (now row 8, was row 7) with more detail.
$ws.Range("I8").Value = "Full geodatabase; May also be able to pull map server with web services: https://hazards.fema.gov/femaportal/wps/portal/NFHLWMS"

# The row insert left the worksheet's Hyperlinks collection anchored to their old
# (pre-insert) addresses instead of following the cells they belonged to, so rebuild
# them all at their correct (post-insert) locations.
while ($ws.Hyperlinks.Count -gt 0) {
    $existing = @($ws.Hyperlinks)
    $existing[0].Delete()
}

# Ozone row link (was H8, now H9).
$ws.Hyperlinks.Add($ws.Range("H9"), "https://ozoneaq.gsfc.nasa.gov/data/omps/", "prods=149") | Out-Null
# Landsat 8 OLI LST row link (unchanged, H4).
$ws.Hyperlinks.Add($ws.Range("H4"), "https://developers.google.com/earth-engine/datasets/catalog/LANDSAT_LC08_C02_T2_L2", "description") | Out-Null
# New Flood Risk web-services link (H8).
$ws.Hyperlinks.Add($ws.Range("H8"), "https://hazards.fema.gov/femaportal/wps/portal/NFHLWMS") | Out-Null

# Keep the hyperlink-styled cells looking like links.
$ws.Range("H4").Style = "Hyperlink"
$ws.Range("H8").Style = "Hyperlink"
$ws.Range("H9").Style = "Hyperlink"

$ws.Range("I9").Select()
